$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "capri black yoga pants",
    "capri bottoms women",
    "capri care",
    "capri compression leggings for women",
    "capri compression pants",
    "capri compression pants men",
    "capri dress leggings",
    "capri for plus size women",
    "capri for women",
    "capri for women sport",
    "capri green leggings for women",
    "capri gym leggings",
    "capri gym pants",
    "capri high waist",
    "capri hiking pants women",
    "capri legging",
    "capri legging black",
    "capri legging for women",
    "capri legging high waist",
    "capri legging pack",
    "capri legging pants",
    "capri legging plus",
    "capri legging workout",
    "capri leggings",
    "capri leggings for summer",
    "capri leggings for women",
    "capri leggings for women black",
    "capri leggings for women pack",
    "capri leggings for women sport",
    "capri leggings for women sport gym",
    "capri leggings for women workout",
    "capri leggings gray",
    "capri leggings high waist",
    "capri leggings medium",
    "capri leggings men",
    "capri leggings pack",
    "capri leggings pack of 6",
    "capri leggings plus",
    "capri leggings sheer",
    "capri leggings skirt",
    "capri leggings small",
    "capri leggings sport",
    "capri leggings under armour",
    "capri leggings with side pocket",
    "capri leggings with skirt",
    "capri leggings women",
    "capri leggings women workout",
    "capri leggings xl",
    "capri leggins",
    "capri length leggings",
    "capri length leggings for women",
    "capri length workout leggings",
    "capri nike",
    "capri pant",
    "capri pants",
    "capri pants 2x",
    "capri pants black",
    "capri pants for women",
    "capri pants women",
    "capri pants yoga",
    "capri pantyhose",
    "capri plus pads",
    "capri plus size leggings",
    "capri pocket leggings",
    "capri running leggings for women",
    "capri running pants",
    "capri running pants for women",
    "capri running pants women",
    "capri running tights",
    "capri running women",
    "capri shaping leggings",
    "capri shorts",
    "capri shorts for women",
    "capri skirt",
    "capri skirt leggings",
    "capri skirts for women",
    "capri sport leggings",
    "capri stockings",
    "capri summer leggings",
    "capri support leggings",
    "capri tight",
    "capri tights",
    "capri tights for women",
    "capri tights for women control top",
    "capri tights men",
    "capri tights women",
    "capri white",
    "capri white leggings",
    "capri with skirt",
    "capri woman",
    "capri womens",
    "capri womens pants",
    "capri workout leggings",
    "capri workout leggings for women",
    "capri workout pants",
    "capri workout pants for women",
    "capri workout pants women",
    "capri xl",
    "capri yoga pants",
    "capri yoga pants for women"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
